$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 1555.6666
$ws.Range("I6").Value = 500.25
$ws.Range("J6").Value = 2400
$ws.Range("K6").Value = 1500.75
$ws.Range("L6").Value = 7200
$ws.Range("M6").Value = -1388.75
$ws.Range("N6").Value = -7424

# Row 8
$ws.Range("H8").Value = 748.6923
$ws.Range("J8").Value = 2966.6667
$ws.Range("L8").Value = 8900.000100000001
$ws.Range("N8").Value = -9178.000100000001

# Row 29
$ws.Range("H29").Value = 3001.5
$ws.Range("I29").Value = 3001.5
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 9004.5
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -8723.5
$ws.Range("N29").ClearContents()

# Row 39
$ws.Range("H39").Value = 780.44446
$ws.Range("I39").Value = 77
$ws.Range("J39").Value = 1659.75
$ws.Range("K39").Value = 231
$ws.Range("L39").Value = 4979.25
$ws.Range("M39").Value = 65
$ws.Range("N39").Value = -5571.25

# Row 62
$ws.Range("H62").Value = 3514
$ws.Range("I62").Value = 2911.111
$ws.Range("J62").Value = 4289.143
$ws.Range("K62").Value = 2911.111
$ws.Range("L62").Value = 4289.143
$ws.Range("M62").Value = -2287.111
$ws.Range("N62").Value = -5537.143

# Row 65
$ws.Range("H65").Value = 3514
$ws.Range("I65").Value = 2911.111
$ws.Range("J65").Value = 4289.143
$ws.Range("K65").Value = 14555.555
$ws.Range("L65").Value = 21445.715
$ws.Range("M65").Value = -11435.555
$ws.Range("N65").Value = -27685.715

$ws = $wb.Worksheets.Item("ARM")
# Row 63
$ws.Range("H63").Value = 4699.9165
$ws.Range("I63").Value = 2342.7144
$ws.Range("J63").Value = 8000
$ws.Range("K63").Value = 2342.7144
$ws.Range("L63").Value = 8000
$ws.Range("M63").Value = -1656.7144
$ws.Range("N63").Value = -9372

# Row 66
$ws.Range("H66").Value = 4699.9165
$ws.Range("I66").Value = 2342.7144
$ws.Range("J66").Value = 8000
$ws.Range("K66").Value = 11713.572
$ws.Range("L66").Value = 40000
$ws.Range("M66").Value = -8281.572
$ws.Range("N66").Value = -46864

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 3611.8125
$ws.Range("I134").Value = 3233.7693
$ws.Range("J134").Value = 5250
$ws.Range("K134").Value = 9701.3079
$ws.Range("L134").Value = 15750
$ws.Range("M134").Value = -7166.3079
$ws.Range("N134").Value = -20820

$ws = $wb.Worksheets.Item("CRP")
# Row 3
$ws.Range("H3").Value = 6712667.5
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 6712667.5
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 6712667.5
$ws.Range("N3").Value = -6712893.5
$ws.Range("M3").ClearContents()

# Row 4
$ws.Range("H4").Value = 70002
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 70002
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 70002
$ws.Range("N4").Value = -70226
$ws.Range("M4").ClearContents()

# Row 17
$ws.Range("H17").Value = 56606
$ws.Range("J17").Value = 56606
$ws.Range("L17").Value = 56606
$ws.Range("N17").Value = -56954

# Row 25
$ws.Range("H25").Value = 47341.668
$ws.Range("I25").Value = 2000
$ws.Range("J25").Value = 70012.5
$ws.Range("K25").Value = 2000
$ws.Range("L25").Value = 70012.5
$ws.Range("M25").Value = -1826
$ws.Range("N25").Value = -70360.5

# Row 35
$ws.Range("H35").Value = 600
$ws.Range("I35").Value = 600
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 600
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -306
$ws.Range("N35").ClearContents()

# Row 59
$ws.Range("H59").Value = 27900
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 27900
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 27900
$ws.Range("N59").Value = -30190
$ws.Range("M59").ClearContents()

# Row 135
$ws.Range("H135").Value = 30000
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 17
$ws.Range("H17").Value = 682.8570999999999
$ws.Range("J17").Value = 733.3333
$ws.Range("L17").Value = 2199.9999
$ws.Range("N17").Value = -2537.9999

# Row 39
$ws.Range("H39").Value = 2911
$ws.Range("J39").Value = 3528.4285
$ws.Range("L39").Value = 10585.2855
$ws.Range("N39").Value = -11173.2855

$ws = $wb.Worksheets.Item("GSM")
# Row 9
$ws.Range("H9").Value = 15088.3
$ws.Range("I9").Value = 844.5
$ws.Range("J9").Value = 36454
$ws.Range("K9").Value = 844.5
$ws.Range("L9").Value = 36454
$ws.Range("M9").Value = -674.5
$ws.Range("N9").Value = -36794

# Row 92
$ws.Range("H92").Value = 14762.5
$ws.Range("J92").Value = 14762.5
$ws.Range("L92").Value = 14762.5
$ws.Range("N92").Value = -18506.5

# Row 113
$ws.Range("H113").Value = 5378.25
$ws.Range("I113").Value = 4000
$ws.Range("J113").Value = 5837.6665
$ws.Range("K113").Value = 4000
$ws.Range("L113").Value = 5837.6665
$ws.Range("M113").Value = -1830
$ws.Range("N113").Value = -10177.6665

$ws = $wb.Worksheets.Item("LTW")
# Row 131
$ws.Range("H131").Value = 18925
$ws.Range("J131").Value = 18925
$ws.Range("L131").Value = 18925
$ws.Range("N131").Value = -29005

# Row 132
$ws.Range("H132").Value = 3844.7778
$ws.Range("I132").Value = 2857.5715
$ws.Range("K132").Value = 8572.7145
$ws.Range("M132").Value = -6042.7145

$ws = $wb.Worksheets.Item("WVR")
# Row 7
$ws.Range("H7").Value = 25102.334
$ws.Range("I7").Value = 201.33333
$ws.Range("J7").Value = 50003.332
$ws.Range("K7").Value = 201.33333
$ws.Range("L7").Value = 50003.332
$ws.Range("M7").Value = -88.33332999999999
$ws.Range("N7").Value = -50229.332

# Row 14
$ws.Range("H14").Value = 38503.5
$ws.Range("I14").Value = 4004
$ws.Range("J14").Value = 50003.332
$ws.Range("K14").Value = 4004
$ws.Range("L14").Value = 50003.332
$ws.Range("N14").Value = -50339.332
$ws.Range("M14").Value = -3836

# Row 132
$ws.Range("H132").Value = 237810.98
$ws.Range("I132").Value = 359263.34
$ws.Range("J132").Value = 11099.866
$ws.Range("K132").Value = 1077790.02
$ws.Range("L132").Value = 33299.598
$ws.Range("M132").Value = -1075260.02
$ws.Range("N132").Value = -38359.598

# Row 136
$ws.Range("H136").Value = 1466.3414
$ws.Range("I136").Value = 703.45715
$ws.Range("K136").Value = 2110.37145
$ws.Range("M136").Value = 439.6285500000004
